{"js": "const body = context.document.body;\nconst ooxml = body.getOoxml();\nawait context.sync();\nreturn ooxml.value.substring(0, 500);\n", "ps1": "# The source diff for this fixture only reorders XML namespace declarations and\n# element attributes (alphabetically) across word/document.xml, the header/footer\n# parts, footnotes.xml and styles.xml -- a canonicalization side-effect of the\n# upstream test-resource tooling being bumped from 2.0.2 to 2.0.3. No paragraph\n# text, run formatting, table content, styles, headers/footers text, or any other\n# observable document content changed. There is no user-visible edit to replay\n# through the Word object model, so we simply touch the document to mirror a\n# no-op re-save and leave all content untouched.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
